$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per the latest crypto data refresh.
# D-column values are forced to remain plain Text (matching the original inlineStr cells)
# so Excel does not auto-convert numeric-looking strings (e.g. "301.16") into numbers,
# and the temporary "@" text format is reverted to the default style afterwards so no
# new cell style is left applied to the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.018.68"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.67%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.602.80"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.82%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("E5").Value = "  +0.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "301.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.12%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3780"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.78%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3629"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.32%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "49.21"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.79%  "

$ws.Range("E10").Value = "  -6.19%  "

$ws.Range("E11").Value = "  +0.14%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08108"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.85%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.79"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.41%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.580"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.18%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.405"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.19%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001243"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.49%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.598.23"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.15%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.15%  "

$ws.Range("E19").Value = "  -1.36%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.59%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.559"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.60%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.5567"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.59%  "

$ws.Range("E23").Value = "  +0.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.48%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "23.011.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.71%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.369"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.95%  "

$ws.Range("E27").Value = "  -4.53%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.07%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "150.41"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.40%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.258"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.14%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "133.74"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.33%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.302"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.32%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.844"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -11.38%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.784.36"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.95%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9615"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.20%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.07630"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.47%  "

$ws.Range("E37").Value = "  -1.75%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.277"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.44%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02712"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.87%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2532"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.39%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.08846"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.82%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.365"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.96%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7048"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.69%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.54"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.07%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6610"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.62%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.001"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.313"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.36%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.990"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.65%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.69"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.00%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07911"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.19%  "

